$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A3").Value = "group3"
$ws2.Range("B3").Value = 0.1
$ws2.Range("C3").Value = 1

$ws2.Range("A4").Select()
$ws2.Activate()
